$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly column: 04_05_2021 -> column AA
$ws.Range("AA1").Value = "04_05_2021"

$ws.Range("AA2").Value  = 12
$ws.Range("AA3").Value  = 16
$ws.Range("AA4").Value  = 33
$ws.Range("AA5").Value  = 47
$ws.Range("AA6").Value  = 113
$ws.Range("AA7").Value  = 271
$ws.Range("AA8").Value  = 402
$ws.Range("AA9").Value  = 562
$ws.Range("AA10").Value = 177
$ws.Range("AA11").Value = 13

# Totals row - extend the running SUM formula into the new column
$ws.Range("AA12").Formula = "=SUM(AA2:AA11)"

# Move selection to reflect where the editor ended up after adding the column
$ws.Range("T11").Select()
